$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Phone"
$ws.Range("B4").Value = 712345678

$ws.Range("A5").Value = "Description"
$ws.Range("B5").Value = "Test RPA"

$ws.Range("C5").Select()
